# Update cryptocurrency price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.754.03'
$ws.Range('D3').Value = '1.738.24'
$ws.Range('E3').Value = '  +5.16%  '
$cell = $ws.Range('D4')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = $origStyle
$ws.Range('E4').Value = '  -0.10%  '
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '227.84'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  +4.16%  '
$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5463'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  +4.28%  '
$ws.Range('E7').Value = '  -0.13%  '
$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.2761'
$cell.Style = $origStyle
$ws.Range('E8').Value = '  +3.21%  '
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06729'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  +5.81%  '
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '22.00'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  +7.34%  '
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.07794'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  +1.47%  '
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.700'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  +2.28%  '
$ws.Range('D13').Value = '1.743.16'
$ws.Range('E13').Value = '  +5.84%  '
$ws.Range('D14').Value = '1.976.52'
$ws.Range('E14').Value = '  +5.10%  '
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5998'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  +6.93%  '
$ws.Range('D16').Value = '0.0₅8434'
$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '69.48'
$cell.Style = $origStyle
$ws.Range('E17').Value = '  +5.90%  '
$ws.Range('D18').Value = '27.753.57'
$ws.Range('E18').Value = '  +6.37%  '
$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '228.44'
$cell.Style = $origStyle
$ws.Range('E19').Value = '  +20.48%  '
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.853'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  +3.69%  '
$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.003'
$cell.Style = $origStyle
$ws.Range('E21').Value = '  -0.12%  '
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.92'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  +6.00%  '
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.242'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  +4.59%  '
$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.004'
$cell.Style = $origStyle
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '147.67'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  +1.25%  '
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.734'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  +13.91%  '
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.1253'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  +4.51%  '
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.472'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  +3.21%  '
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.22'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  +8.17%  '
$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.05701'
$cell.Style = $origStyle
$ws.Range('E30').Value = '  +1.33%  '
$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.315'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  +3.59%  '
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.706'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  +6.29%  '
$ws.Range('E33').Value = '  +4.36%  '
$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.693'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  +7.19%  '
$cell = $ws.Range('D35')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.9768'
$cell.Style = $origStyle
$ws.Range('E35').Value = '  +3.26%  '
$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.856'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  +2.08%  '
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.447'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  +1.49%  '
$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5995'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  +4.17%  '
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.01673'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  +5.25%  '
$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.946'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('D41').Value = '1.052.34'
$ws.Range('E41').Value = '  +2.92%  '
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.8493'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('E43').Value = '  -0.06%  '
$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '102.15'
$cell.Style = $origStyle
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('D45').Value = '1.880.65'
$ws.Range('E45').Value = '  +5.01%  '
$ws.Range('E46').Value = '  +13.47%  '
$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '59.69'
$cell.Style = $origStyle
$ws.Range('E47').Value = '  +2.43%  '
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.308'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  +3.82%  '
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.4437'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  +2.13%  '
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.004'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  +0.11%  '
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.05334'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  +0.73%  '
